$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 6.206015333333333
$ws.Range("H2").Value2 = 18.618046
$ws.Range("I2").Value2 = 0.0150172404156507
$ws.Range("J2").Value2 = 0.0150172404156507
$ws.Range("M2").Value2 = 61.04160633333334
$ws.Range("N2").Value2 = 183.124819
$ws.Range("O2").Value2 = 0.2043613460574534
$ws.Range("P2").Value2 = 0.2043613460574534
$ws.Range("Q2").Value2 = 378.8251448759638
$ws.Range("R2").Value2 = 3409.426303883674
$ws.Range("S2").Value2 = 0.003068943465410768
$ws.Range("T2").Value2 = 0.003068943465410768

$ws.Range("G3").Value2 = 6.206015333333333
$ws.Range("H3").Value2 = 18.618046
$ws.Range("I3").Value2 = 0.0150172404156507
$ws.Range("J3").Value2 = 0.0150172404156507
$ws.Range("O3").Value2 = 0.3559304658284363
$ws.Range("P3").Value2 = 0.3559304658284363
$ws.Range("Q3").Value2 = 659.7892061511453
$ws.Range("R3").Value2 = 5938.102855360308
$ws.Range("S3").Value2 = 0.005345093376600173
$ws.Range("T3").Value2 = 0.005345093376600174

$ws.Range("G4").Value2 = 6.206015333333333
$ws.Range("H4").Value2 = 18.618046
$ws.Range("I4").Value2 = 0.0150172404156507
$ws.Range("J4").Value2 = 0.0150172404156507
$ws.Range("M4").Value2 = 131.3384093333333
$ws.Range("N4").Value2 = 394.015228
$ws.Range("O4").Value2 = 0.4397081881141102
$ws.Range("P4").Value2 = 0.4397081881141103
$ws.Range("Q4").Value2 = 815.0881821782764
$ws.Range("R4").Value2 = 7335.793639604488
$ws.Range("S4").Value2 = 0.006603203573639757
$ws.Range("T4").Value2 = 0.006603203573639757

$ws.Range("I5").Value2 = 0.9317452840597572
$ws.Range("J5").Value2 = 0.9317452840597571
$ws.Range("M5").Value2 = 61.04160633333334
$ws.Range("N5").Value2 = 183.124819
$ws.Range("O5").Value2 = 0.2043613460574534
$ws.Range("P5").Value2 = 0.2043613460574534
$ws.Range("Q5").Value2 = 23504.22131176485
$ws.Range("R5").Value2 = 211537.9918058837
$ws.Range("S5").Value2 = 0.1904127204331363
$ws.Range("T5").Value2 = 0.1904127204331363

$ws.Range("I6").Value2 = 0.9317452840597572
$ws.Range("J6").Value2 = 0.9317452840597571
$ws.Range("O6").Value2 = 0.3559304658284363
$ws.Range("P6").Value2 = 0.3559304658284363
$ws.Range("S6").Value2 = 0.3316365329888381
$ws.Range("T6").Value2 = 0.3316365329888381

$ws.Range("I7").Value2 = 0.9317452840597572
$ws.Range("J7").Value2 = 0.9317452840597571
$ws.Range("M7").Value2 = 131.3384093333333
$ws.Range("N7").Value2 = 394.015228
$ws.Range("O7").Value2 = 0.4397081881141102
$ws.Range("P7").Value2 = 0.4397081881141103
$ws.Range("Q7").Value2 = 50572.17896345052
$ws.Range("R7").Value2 = 455149.6106710547
$ws.Range("S7").Value2 = 0.4096960306377828
$ws.Range("T7").Value2 = 0.4096960306377828

$ws.Range("G8").Value2 = 22.00088566666667
$ws.Range("H8").Value2 = 66.002657
$ws.Range("I8").Value2 = 0.05323747552459213
$ws.Range("J8").Value2 = 0.05323747552459213
$ws.Range("M8").Value2 = 61.04160633333334
$ws.Range("N8").Value2 = 183.124819
$ws.Range("O8").Value2 = 0.2043613460574534
$ws.Range("P8").Value2 = 0.2043613460574534
$ws.Range("Q8").Value2 = 1342.969401849343
$ws.Range("R8").Value2 = 12086.72461664409
$ws.Range("S8").Value2 = 0.01087968215890638
$ws.Range("T8").Value2 = 0.01087968215890638

$ws.Range("G9").Value2 = 22.00088566666667
$ws.Range("H9").Value2 = 66.002657
$ws.Range("I9").Value2 = 0.05323747552459213
$ws.Range("J9").Value2 = 0.05323747552459213
$ws.Range("O9").Value2 = 0.3559304658284363
$ws.Range("P9").Value2 = 0.3559304658284363
$ws.Range("Q9").Value2 = 2339.01241117872
$ws.Range("R9").Value2 = 21051.11170060849
$ws.Range("S9").Value2 = 0.01894883946299806
$ws.Range("T9").Value2 = 0.01894883946299806

$ws.Range("G10").Value2 = 22.00088566666667
$ws.Range("H10").Value2 = 66.002657
$ws.Range("I10").Value2 = 0.05323747552459213
$ws.Range("J10").Value2 = 0.05323747552459213
$ws.Range("M10").Value2 = 131.3384093333333
$ws.Range("N10").Value2 = 394.015228
$ws.Range("O10").Value2 = 0.4397081881141102
$ws.Range("P10").Value2 = 0.4397081881141103
$ws.Range("Q10").Value2 = 2889.561327384532
$ws.Range("R10").Value2 = 26006.05194646079
$ws.Range("S10").Value2 = 0.0234089539026877
$ws.Range("T10").Value2 = 0.0234089539026877
